# Apply updated PnL risk-driver calibration estimates (std moving average)
$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.2179356199638502
$wsLinear.Range("B3").Value = 1.084517857588632
$wsLinear.Range("B4").Value = 1335.125292586363
$wsLinear.Range("B5").Value = "[1.0, 0.12659801392715944, 0.012964859385988294, -0.024861606153344415, -0.037538340786715994, -0.020849437985720507, 0.13925338354901312, 0.2909496510151006, 0.11467813181256743, -0.03199967872744598, -0.05021818161594246, -0.045471915513467645, -0.03569855726553589, 0.12004018162210706, 0.31108107334952145, 0.14524462871693625, -0.028176057717056474, -0.04187974157321775, -0.03193568637929988, -0.03756239412689398]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 1.126529954553905
$wsNonLinear.Range("B5").Value = 4.077400978062738
$wsNonLinear.Range("B6").Value = 1286.404987171449
$wsNonLinear.Range("B7").Value = 2.946239134196051
$wsNonLinear.Range("B8").Value = -4.537919863065801
$wsNonLinear.Range("B9").Value = 1375.932037524403
$wsNonLinear.Range("B10").Value = "[0.9999999999999999, 0.12403499288319315, 0.015721475562531267, -0.020705187569092696, -0.033485573719376265, -0.017226652452552022, 0.1376975416491414, 0.2832238307515278, 0.11357661601956445, -0.026567293889660733, -0.04630152829326754, -0.04120205061633143, -0.03081530388543582, 0.11846509740561778, 0.3043226804168205, 0.14259217338675245, -0.022749345617106113, -0.03869660766932291, -0.026184827016304174, -0.032181740388883426]"
